# Remove the redundant "труб" (pipe/tube) word from the duration-by-TCP
# stepwise-extrapolation-descending template.
#
# Before: "...Общая протяженность проектируемой сети из %PM% труб – %PL% км."
# After:  "...Общая протяженность проектируемой сети из %PM% – %PL% км."
#
# Editing the text at that point is exactly what a user would do with
# Word's Find & Replace (Ctrl+H): Word automatically drops/updates the
# "_GoBack" bookmark at the location of the last edit, which is why the
# canonical OOXML diff for this change also shows a new _GoBack bookmark
# appearing right where "труб " used to be (and the pre-existing
# "_Hlk93481827" bookmark's id shifting up by one to make room for it).

$d = $word.ActiveDocument

# Find the exact run of text to remove: the word "труб" plus the single
# trailing space that follows it (the dash "–" that follows is kept).
$found = $d.Content.Duplicate
[void]$found.Find.Execute("труб ", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)

if ($found.Find.Found) {
    # Mark the edit point with a "_GoBack" bookmark collapsed at the start
    # of the text being removed -- mirrors Word's own automatic behavior.
    $editPoint = $d.Range($found.Start, $found.Start)
    $d.Bookmarks.Add("_GoBack", $editPoint)

    # Now actually delete "труб " so the sentence reads "...%PM% – %PL%...".
    $found.Delete()
}
